$d = $word.ActiveDocument

# --- Locate the "Heading 3" paragraph and position right after it ---
$findRange = $d.Content
$found = $findRange.Find.Execute("Heading 3", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headingPara = $findRange.Paragraphs(1)
$insertPos = $headingPara.Range.End
$rng = $d.Range($insertPos, $insertPos)

# --- Insert the new BodyText paragraph + filerequirements table + trailing BodyText paragraph ---
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p><w:tbl><w:tblPr><w:tblW w:w="5000" w:type="pct"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="0020" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="1213"/><w:gridCol w:w="5121"/><w:gridCol w:w="781"/><w:gridCol w:w="781"/><w:gridCol w:w="820"/><w:gridCol w:w="860"/></w:tblGrid><w:tr><w:trPr><w:tblHeader/></w:trPr><w:tc><w:tcPr><w:tcW w:w="1002" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Formatted file</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4234" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Description</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>1. QC screening</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>2. QC reporting</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="678" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>3. Data analysis</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="711" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>4. WQX formatting</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1002" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>Results</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4234" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>Water quality results organized by sample location and date</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="678" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="711" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1002" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>DQO accuracy</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4234" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t xml:space="preserve">Summary of data quality objectives that describe quality control </w:t></w:r><w:r><w:t>accuracy</w:t></w:r><w:r><w:t xml:space="preserve"> for data in the results file</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="711" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1002" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>DQO frequency and completeness</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4234" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>Summary of data quality objectives that describe quality control frequency and completeness measures for data in the results file</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="711" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1002" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>Sites</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4234" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>A site metadata file, including location names, latitude, longitude, and additional grouping factors for sites in the results file</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="678" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="711" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1002" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>WQX metadata</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4234" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>wqx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> metadata file required for generating output to facilitate data upload to WQX</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="646" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="678" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="711" w:type="dxa"/><w:tcBorders><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders></w:tcPr><w:p><w:pPr><w:pStyle w:val="Compact"/></w:pPr><w:r><w:t>&#x2713;</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)

# --- Update the "Compact" paragraph style: single line spacing + explicit 10pt run size ---
$compact = $d.Styles("Compact")
$compact.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle
$compact.Font.Size = 10

Write-Host ("Paragraphs: " + $d.Paragraphs.Count)
Write-Host ("Tables: " + $d.Tables.Count)
